$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert a new row at row 3 (Prejuveniles caballeros gains a new #1: Petric, Juan Cruz)
$ws.Rows.Item(3).Insert()

# Insert a new row at row 10 (after the row-3 insert shifts everything down by one,
# the Sub23 caballeros block now starts at row 10; insert a new #1 there: Suarez, Fermin Iru)
$ws.Rows.Item(10).Insert()

$torneo = "####1er Torneo Federativo - C.A.E. - Sub 23, Prejuveniles y sub 23 (28 de Febrero y 1 de Marzo) - Juniors (Domingo 1 de Marzo)"

# New row 3: Petric, Juan Cruz (Prejuveniles / caballeros / 1)
$ws.Cells.Item(3, 1).Value = $torneo
$ws.Cells.Item(3, 2).Value = "Prejuveniles"
$ws.Cells.Item(3, 3).Value = "caballeros"
$ws.Cells.Item(3, 4).Value = 1
$ws.Cells.Item(3, 5).Value = "Petric, Juan Cruz"
$ws.Cells.Item(3, 6).Value = 81
$ws.Cells.Item(3, 7).Value = ""
$ws.Cells.Item(3, 8).Value = 81

# Existing Prejuveniles caballeros rows shift down one and their posicion increments by 1
$ws.Cells.Item(4, 4).Value = 2
$ws.Cells.Item(5, 4).Value = 3
$ws.Cells.Item(6, 4).Value = 4

# New row 10: Suarez, Fermin Iru (Sub23 / caballeros / 1)
$ws.Cells.Item(10, 1).Value = $torneo
$ws.Cells.Item(10, 2).Value = "Sub23"
$ws.Cells.Item(10, 3).Value = "caballeros"
$ws.Cells.Item(10, 4).Value = 1
$ws.Cells.Item(10, 5).Value = "Suárez, Fermín Irú"
$ws.Cells.Item(10, 6).Value = 78
$ws.Cells.Item(10, 7).Value = ""
$ws.Cells.Item(10, 8).Value = 78

# Existing Sub23 caballeros rows shift down one and their posicion increments by 1
$ws.Cells.Item(11, 4).Value = 2
$ws.Cells.Item(12, 4).Value = 3

$ws.Range("A1:H12").EntireColumn.AutoFit() | Out-Null
